$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Citation 1: "...and Broschat, S.L., Submitted. PARGT: A Software Tool..."
#    becomes "...and Broschat, S.L., 2020. PARGT: A Software Tool..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Submitted. PARGT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2020. PARGT", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the stray "_GoBack" bookmark that sat right after "Scientific
#    Reports" in that same citation.
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    if ($goBack -ne $null) {
        $goBack.Delete()
    }
} catch {
    # no _GoBack bookmark present - nothing to do
}

# ---------------------------------------------------------------------------
# 3) Append a new third citation entry to the numbered "Citations" list,
#    reusing the same list style/numbering as the two existing entries.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "Chowdhury, A.S., Khaledian, E. and Broschat, S.L., 2019. " + `
    "Capreomycin resistance prediction in two species of Mycobacterium " + `
    "using a stacked ensemble method. Journal of applied microbiology."

# Italicize the genus name "Mycobacterium" within the new citation.
$rng = $newPara.Range.Duplicate
$found = $rng.Find.Execute("Mycobacterium", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Italic = $true
}
